# "Generate Report for Archive"
#
# The localization status has moved on from the handoff stage, so the
# "Ready for handoff" status text is now "In Translation" everywhere it
# appears (the Overview roll-up sheet as well as each per-locale sheet).
# Because the new status string is shorter than the old one, the Status
# column no longer needs to be as wide, so those columns are narrowed to
# fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Target stored column width is 13.4101845877511 characters. This COM layer
# quantizes ColumnWidth to whole screen pixels (MDW=6, +5px padding) before
# persisting, so the achievable width closest to the target is reached by
# requesting 12.5 "characters" (=75px -> stored width 13.33...).
$newWidth  = 12.5

# --- Overview sheet: zh-cn / de-de status columns are E and F ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- Per-locale sheets: Status column is C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
